$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3134.7102
$ws.Range("J17").Value = 3205.8955
$ws.Range("L17").Value = 9617.6865
$ws.Range("N17").Value = -9953.6865

$ws.Range("H113").Value = 4318.9287
$ws.Range("I113").Value = 3423.5715
$ws.Range("J113").Value = 5214.2856
$ws.Range("K113").Value = 3423.5715
$ws.Range("L113").Value = 5214.2856
$ws.Range("M113").Value = -169.5715
$ws.Range("N113").Value = -11722.2856

$ws.Range("H116").Value = 2368.0908
$ws.Range("I116").Value = 2686.625
$ws.Range("K116").Value = 2686.625
$ws.Range("M116").Value = 755.375

$ws.Range("H129").Value = 1264.7
$ws.Range("I129").Value = 291.16666
$ws.Range("J129").Value = 2725
$ws.Range("K129").Value = 873.4999799999999
$ws.Range("L129").Value = 8175
$ws.Range("M129").Value = 4126.50002
$ws.Range("N129").Value = -18175

$ws.Range("H132").Value = 26069.334
$ws.Range("I132").Value = 26680.879
$ws.Range("K132").Value = 80042.637
$ws.Range("M132").Value = -77512.637

$ws.Range("H135").Value = 3135.6
$ws.Range("I135").Value = 2876
$ws.Range("J135").Value = 3395.2
$ws.Range("K135").Value = 25884
$ws.Range("L135").Value = 30556.8
$ws.Range("M135").Value = -23349
$ws.Range("N135").Value = -35626.8

$ws.Range("H137").Value = 37039468
$ws.Range("I137").Value = 55556916
$ws.Range("J137").Value = 4577.778
$ws.Range("K137").Value = 166670748
$ws.Range("L137").Value = 13733.334
$ws.Range("M137").Value = -166668198
$ws.Range("N137").Value = -18833.334

$ws.Range("H138").Value = 5726904.5
$ws.Range("I138").Value = 2085991.9
$ws.Range("J138").Value = 7465847.5
$ws.Range("K138").Value = 6257975.699999999
$ws.Range("L138").Value = 22397542.5
$ws.Range("M138").Value = -6252835.699999999
$ws.Range("N138").Value = -22407822.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21549
$ws.Range("I32").Value = 6197.7163
$ws.Range("K32").Value = 6197.7163
$ws.Range("M32").Value = -5910.7163

$ws.Range("H61").Value = 2841.4285
$ws.Range("I61").Value = 2319.9412
$ws.Range("J61").Value = 5057.75
$ws.Range("K61").Value = 2319.9412
$ws.Range("L61").Value = 5057.75
$ws.Range("M61").Value = -2107.9412
$ws.Range("N61").Value = -5481.75

$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 2239.4119
$ws.Range("I132").Value = 1629.6072
$ws.Range("J132").Value = 5085.1665
$ws.Range("K132").Value = 4888.821599999999
$ws.Range("L132").Value = 15255.4995
$ws.Range("M132").Value = -2358.821599999999
$ws.Range("N132").Value = -20315.4995

$ws.Range("H133").Value = 47250
$ws.Range("J133").Value = 47250
$ws.Range("L133").Value = 47250
$ws.Range("N133").Value = -52310

$ws.Range("H136").Value = 2841.4285
$ws.Range("I136").Value = 2319.9412
$ws.Range("J136").Value = 5057.75
$ws.Range("K136").Value = 6959.823600000001
$ws.Range("L136").Value = 15173.25
$ws.Range("M136").Value = -4409.823600000001
$ws.Range("N136").Value = -20273.25

$ws.Range("H139").Value = 47107.375
$ws.Range("J139").Value = 47107.375
$ws.Range("L139").Value = 47107.375
$ws.Range("N139").Value = -57387.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1637.25
$ws.Range("I20").Value = 1469.6
$ws.Range("J20").Value = 1916.6666
$ws.Range("K20").Value = 1469.6
$ws.Range("L20").Value = 1916.6666
$ws.Range("M20").Value = -1222.6
$ws.Range("N20").Value = -2410.6666

$ws.Range("H99").Value = 1928.8182
$ws.Range("I99").Value = 1928.8182
$ws.Range("K99").Value = 1928.8182
$ws.Range("M99").Value = -430.8181999999999

$ws.Range("H105").Value = 3779.1667
$ws.Range("I105").Value = 3687.5
$ws.Range("J105").Value = 3962.5
$ws.Range("K105").Value = 3687.5
$ws.Range("L105").Value = 3962.5
$ws.Range("M105").Value = -1940.5
$ws.Range("N105").Value = -7456.5

$ws.Range("H107").Value = 2701.6667
$ws.Range("I107").Value = 2801
$ws.Range("J107").Value = 2304.3333
$ws.Range("K107").Value = 2801
$ws.Range("L107").Value = 2304.3333
$ws.Range("M107").Value = -881
$ws.Range("N107").Value = -6144.3333

$ws.Range("H134").Value = 2005.6078
$ws.Range("I134").Value = 1663.591
$ws.Range("J134").Value = 4155.4287
$ws.Range("K134").Value = 4990.772999999999
$ws.Range("L134").Value = 12466.2861
$ws.Range("M134").Value = -2455.772999999999
$ws.Range("N134").Value = -17536.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 808.3
$ws.Range("I22").Value = 564.7778
$ws.Range("K22").Value = 564.7778
$ws.Range("M22").Value = -214.7778

$ws.Range("H31").Value = 3993.8108
$ws.Range("I31").Value = 1889.2059
$ws.Range("J31").Value = 5782.725
$ws.Range("K31").Value = 1889.2059
$ws.Range("L31").Value = 5782.725
$ws.Range("M31").Value = -1594.2059
$ws.Range("N31").Value = -6372.725

$ws.Range("H34").Value = 3993.8108
$ws.Range("I34").Value = 1889.2059
$ws.Range("J34").Value = 5782.725
$ws.Range("K34").Value = 1889.2059
$ws.Range("L34").Value = 5782.725
$ws.Range("M34").Value = -1687.2059
$ws.Range("N34").Value = -6186.725

$ws.Range("H44").Value = 5300
$ws.Range("I44").Value = 2000
$ws.Range("K44").Value = 2000
$ws.Range("M44").Value = -1558

$ws.Range("H55").Value = 4700
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 7400
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 7400
$ws.Range("M55").Value = -1685
$ws.Range("N55").Value = -8030

$ws.Range("H99").Value = 15154228
$ws.Range("I99").Value = 2430.9333
$ws.Range("J99").Value = 47622364
$ws.Range("K99").Value = 2430.9333
$ws.Range("L99").Value = 47622364
$ws.Range("M99").Value = -932.9333000000001
$ws.Range("N99").Value = -47625360

$ws.Range("H102").Value = 26466.666
$ws.Range("I102").Value = 20000
$ws.Range("J102").Value = 29700
$ws.Range("K102").Value = 20000
$ws.Range("L102").Value = 29700
$ws.Range("M102").Value = -17566
$ws.Range("N102").Value = -34568

$ws.Range("H126").Value = 15154228
$ws.Range("I126").Value = 2430.9333
$ws.Range("J126").Value = 47622364
$ws.Range("K126").Value = 7292.7999
$ws.Range("L126").Value = 142867092
$ws.Range("M126").Value = -4822.7999
$ws.Range("N126").Value = -142872032

$ws.Range("H134").Value = 30614778
$ws.Range("I134").Value = 38463384
$ws.Range("J134").Value = 21742438
$ws.Range("K134").Value = 115390152
$ws.Range("L134").Value = 65227314
$ws.Range("M134").Value = -115387617
$ws.Range("N134").Value = -65232384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 830.1875
$ws.Range("I23").Value = 2593.25
$ws.Range("J23").Value = 242.5
$ws.Range("K23").Value = 7779.75
$ws.Range("L23").Value = 727.5
$ws.Range("M23").Value = -7544.75
$ws.Range("N23").Value = -1197.5

$ws.Range("H107").Value = 545.5
$ws.Range("I107").Value = 573.3077
$ws.Range("J107").Value = 517.6923
$ws.Range("K107").Value = 1719.9231
$ws.Range("L107").Value = 1553.0769
$ws.Range("M107").Value = 200.0769
$ws.Range("N107").Value = -5393.0769

$ws.Range("H132").Value = 1631.7778
$ws.Range("I132").Value = 1371.5
$ws.Range("K132").Value = 12343.5
$ws.Range("M132").Value = -9813.5

$ws.Range("H136").Value = 2758.7693
$ws.Range("I136").Value = 1743.3334
$ws.Range("J136").Value = 2862.034
$ws.Range("K136").Value = 5230.0002
$ws.Range("L136").Value = 8586.102000000001
$ws.Range("M136").Value = -130.0002000000004
$ws.Range("N136").Value = -18786.102

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4973.5713
$ws.Range("I102").Value = 3224.6428
$ws.Range("J102").Value = 8471.429
$ws.Range("K102").Value = 3224.6428
$ws.Range("L102").Value = 8471.429
$ws.Range("M102").Value = -1602.6428
$ws.Range("N102").Value = -11715.429

$ws.Range("H122").Value = 2978.9285
$ws.Range("I122").Value = 2977.3076
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8931.9228
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6481.9228
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 3056.4102
$ws.Range("I126").Value = 3183.3333
$ws.Range("K126").Value = 9549.999899999999
$ws.Range("M126").Value = -7079.999899999999

$ws.Range("H138").Value = 89300
$ws.Range("J138").Value = 89300
$ws.Range("L138").Value = 89300
$ws.Range("N138").Value = -99580

$ws.Range("H139").Value = 45863
$ws.Range("J139").Value = 45863
$ws.Range("L139").Value = 45863
$ws.Range("N139").Value = -56143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4834.3335
$ws.Range("I93").Value = 4003
$ws.Range("J93").Value = 5250
$ws.Range("K93").Value = 4003
$ws.Range("L93").Value = 5250
$ws.Range("M93").Value = -2755
$ws.Range("N93").Value = -7746

$ws.Range("H122").Value = 3513.8147
$ws.Range("I122").Value = 2690.3
$ws.Range("K122").Value = 8070.900000000001
$ws.Range("M122").Value = -5620.900000000001

$ws.Range("H132").Value = 3136.5122
$ws.Range("I132").Value = 2223.6667
$ws.Range("K132").Value = 6671.000100000001
$ws.Range("M132").Value = -4141.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1159.2826
$ws.Range("I122").Value = 1090.9487
$ws.Range("K122").Value = 3272.8461
$ws.Range("M122").Value = -822.8460999999998

$ws.Range("H126").Value = 48394.43
$ws.Range("I126").Value = 55979.945
$ws.Range("K126").Value = 167939.835
$ws.Range("M126").Value = -165469.835

$ws.Range("H132").Value = 3365.9185
$ws.Range("I132").Value = 3070.182
$ws.Range("K132").Value = 9210.545999999998
$ws.Range("M132").Value = -6680.545999999998
